$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.632000000000001
$ws.Range("C10").Value = -12.937
$ws.Range("C12").Value = -11.066
$ws.Range("D15").Value = -8.204000000000001
$ws.Range("C18").Value = -12.843
$ws.Range("D20").Value = -7.517000000000001
$ws.Range("D29").Value = -7.292
$ws.Range("D30").Value = -7.178999999999999
$ws.Range("D31").Value = -8.023
$ws.Range("C37").Value = -13.243
$ws.Range("D40").Value = -7.641
$ws.Range("C55").Value = -13.916
$ws.Range("C68").Value = -11.167
$ws.Range("D68").Value = -6.879
$ws.Range("D76").Value = -7.311999999999999
$ws.Range("C77").Value = -13.117
$ws.Range("C78").Value = -13.214
$ws.Range("D87").Value = -8.305
$ws.Range("D88").Value = -7.709000000000001
$ws.Range("D96").Value = -7.267
$ws.Range("D98").Value = -8.228
$ws.Range("D101").Value = -7.886000000000001
$ws.Range("D102").Value = -8.036
